# Auto-generated Excel COM-interop script
# Applies the numeric value updates described by the commit diff
# to the "Leve" profit-tracking tables across all 8 class sheets.
$wb = $excel.ActiveWorkbook

# --- ALC!row 15 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 798.37805
$ws.Cells.Item(15, 9).Value = 798.37805
$ws.Cells.Item(15, 11).Value = 2395.13415
$ws.Cells.Item(15, 13).Value = -2226.13415

# --- ALC!row 92 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(92, 8).Value = 17857416
$ws.Cells.Item(92, 9).Value = 19230986
$ws.Cells.Item(92, 10).Value = 998
$ws.Cells.Item(92, 11).Value = 19230986
$ws.Cells.Item(92, 12).Value = 998
$ws.Cells.Item(92, 13).Value = -19229738
$ws.Cells.Item(92, 14).Value = -3494

# --- ALC!row 98 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2200.1462
$ws.Cells.Item(98, 9).Value = 1902.8064
$ws.Cells.Item(98, 10).Value = 3121.9
$ws.Cells.Item(98, 11).Value = 1902.8064
$ws.Cells.Item(98, 12).Value = 3121.9
$ws.Cells.Item(98, 13).Value = -404.8063999999999
$ws.Cells.Item(98, 14).Value = -6117.9

# --- ALC!row 103 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(103, 8).Value = 924.34784
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).ClearContents()

# --- ALC!row 107 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(107, 8).Value = 1091.5
$ws.Cells.Item(107, 10).Value = 1233.1666
$ws.Cells.Item(107, 12).Value = 1233.1666
$ws.Cells.Item(107, 14).Value = -5073.1666

# --- ALC!row 122 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 2200.1462
$ws.Cells.Item(122, 9).Value = 1902.8064
$ws.Cells.Item(122, 10).Value = 3121.9
$ws.Cells.Item(122, 11).Value = 5708.4192
$ws.Cells.Item(122, 12).Value = 9365.700000000001
$ws.Cells.Item(122, 13).Value = -3258.4192
$ws.Cells.Item(122, 14).Value = -14265.7

# --- ALC!row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 92545.09
$ws.Cells.Item(137, 9).Value = 2001
$ws.Cells.Item(137, 11).Value = 6003
$ws.Cells.Item(137, 13).Value = -3453

# --- ARM!row 2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 231809.58
$ws.Cells.Item(2, 9).Value = 308933
$ws.Cells.Item(2, 10).Value = 439.33334
$ws.Cells.Item(2, 11).Value = 308933
$ws.Cells.Item(2, 12).Value = 439.33334
$ws.Cells.Item(2, 13).Value = -308820
$ws.Cells.Item(2, 14).Value = -665.33334

# --- ARM!row 22 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(22, 8).Value = 1000000
$ws.Cells.Item(22, 9).Value = 1000000
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 1000000
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -999701
$ws.Cells.Item(22, 14).ClearContents()

# --- ARM!row 32 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3551.75
$ws.Cells.Item(32, 9).Value = 2117.221
$ws.Cells.Item(32, 10).Value = 15888.7
$ws.Cells.Item(32, 11).Value = 2117.221
$ws.Cells.Item(32, 12).Value = 15888.7
$ws.Cells.Item(32, 13).Value = -1830.221
$ws.Cells.Item(32, 14).Value = -16462.7

# --- ARM!row 45 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1436.8235
$ws.Cells.Item(45, 9).Value = 1103.7778
$ws.Cells.Item(45, 11).Value = 1103.7778
$ws.Cells.Item(45, 13).Value = -726.7778000000001

# --- ARM!row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1406.4546
$ws.Cells.Item(102, 9).Value = 1406.4546
$ws.Cells.Item(102, 11).Value = 1406.4546
$ws.Cells.Item(102, 13).Value = 215.5454

# --- ARM!row 116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(116, 8).Value = 231809.58
$ws.Cells.Item(116, 9).Value = 308933
$ws.Cells.Item(116, 10).Value = 439.33334
$ws.Cells.Item(116, 11).Value = 308933
$ws.Cells.Item(116, 12).Value = 439.33334
$ws.Cells.Item(116, 13).Value = -306639
$ws.Cells.Item(116, 14).Value = -5027.33334

# --- ARM!row 122 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 1491.5264
$ws.Cells.Item(122, 9).Value = 1019.94116
$ws.Cells.Item(122, 10).Value = 5500
$ws.Cells.Item(122, 11).Value = 3059.82348
$ws.Cells.Item(122, 12).Value = 16500
$ws.Cells.Item(122, 13).Value = -609.82348
$ws.Cells.Item(122, 14).Value = -21400

# --- ARM!row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 3530.6365
$ws.Cells.Item(132, 9).Value = 2898.625
$ws.Cells.Item(132, 11).Value = 8695.875
$ws.Cells.Item(132, 13).Value = -6165.875

# --- BSM!row 3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 231809.58
$ws.Cells.Item(3, 9).Value = 308933
$ws.Cells.Item(3, 10).Value = 439.33334
$ws.Cells.Item(3, 11).Value = 308933
$ws.Cells.Item(3, 12).Value = 439.33334
$ws.Cells.Item(3, 13).Value = -308819
$ws.Cells.Item(3, 14).Value = -667.33334

# --- BSM!row 94 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 785.5
$ws.Cells.Item(94, 9).Value = 338.375
$ws.Cells.Item(94, 10).Value = 2574
$ws.Cells.Item(94, 11).Value = 338.375
$ws.Cells.Item(94, 12).Value = 2574
$ws.Cells.Item(94, 13).Value = 112.625
$ws.Cells.Item(94, 14).Value = -3476

# --- BSM!row 105 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2268
$ws.Cells.Item(105, 9).Value = 2093
$ws.Cells.Item(105, 10).Value = 3213
$ws.Cells.Item(105, 11).Value = 2093
$ws.Cells.Item(105, 12).Value = 3213
$ws.Cells.Item(105, 13).Value = -346
$ws.Cells.Item(105, 14).Value = -6707

# --- BSM!row 107 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1071.6666
$ws.Cells.Item(107, 9).Value = 1106
$ws.Cells.Item(107, 11).Value = 1106
$ws.Cells.Item(107, 13).Value = 814

# --- BSM!row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 6099.5356
$ws.Cells.Item(134, 9).Value = 7973.8423
$ws.Cells.Item(134, 10).Value = 2142.6667
$ws.Cells.Item(134, 11).Value = 23921.5269
$ws.Cells.Item(134, 12).Value = 6428.000100000001
$ws.Cells.Item(134, 13).Value = -21386.5269
$ws.Cells.Item(134, 14).Value = -11498.0001

# --- CRP!row 16 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 938.75
$ws.Cells.Item(16, 9).Value = 930
$ws.Cells.Item(16, 11).Value = 930
$ws.Cells.Item(16, 13).Value = -643

# --- CRP!row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2591.5334
$ws.Cells.Item(31, 9).Value = 1754.5714
$ws.Cells.Item(31, 11).Value = 1754.5714
$ws.Cells.Item(31, 13).Value = -1459.5714

# --- CRP!row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2591.5334
$ws.Cells.Item(34, 9).Value = 1754.5714
$ws.Cells.Item(34, 11).Value = 1754.5714
$ws.Cells.Item(34, 13).Value = -1552.5714

# --- CRP!row 113 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 938.75
$ws.Cells.Item(113, 9).Value = 930
$ws.Cells.Item(113, 11).Value = 930
$ws.Cells.Item(113, 13).Value = 1240

# --- CRP!row 122 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2258.4666
$ws.Cells.Item(122, 9).Value = 1145.3
$ws.Cells.Item(122, 11).Value = 3435.9
$ws.Cells.Item(122, 13).Value = -985.8999999999996

# --- CRP!row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 2158.5833
$ws.Cells.Item(132, 9).Value = 1590.5
$ws.Cells.Item(132, 11).Value = 4771.5
$ws.Cells.Item(132, 13).Value = -2241.5

# --- CRP!row 134 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2428.2173
$ws.Cells.Item(134, 9).Value = 2255.8667
$ws.Cells.Item(134, 11).Value = 6767.6001
$ws.Cells.Item(134, 13).Value = -4232.6001

# --- CUL!row 4 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1634542.8
$ws.Cells.Item(4, 9).Value = 2050108.8
$ws.Cells.Item(4, 11).Value = 6150326.4
$ws.Cells.Item(4, 13).Value = -6150214.4

# --- CUL!row 5 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 560.9231
$ws.Cells.Item(5, 9).Value = 510.45
$ws.Cells.Item(5, 10).Value = 729.1667
$ws.Cells.Item(5, 11).Value = 1531.35
$ws.Cells.Item(5, 12).Value = 2187.5001
$ws.Cells.Item(5, 13).Value = -1419.35
$ws.Cells.Item(5, 14).Value = -2411.5001

# --- CUL!row 68 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 883.8333
$ws.Cells.Item(68, 9).Value = 680
$ws.Cells.Item(68, 10).Value = 902.36365
$ws.Cells.Item(68, 11).Value = 2040
$ws.Cells.Item(68, 12).Value = 2707.09095
$ws.Cells.Item(68, 13).Value = -1229
$ws.Cells.Item(68, 14).Value = -4329.09095

# --- CUL!row 71 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(71, 8).Value = 883.8333
$ws.Cells.Item(71, 9).Value = 680
$ws.Cells.Item(71, 10).Value = 902.36365
$ws.Cells.Item(71, 11).Value = 6120
$ws.Cells.Item(71, 12).Value = 8121.27285
$ws.Cells.Item(71, 13).Value = -2064
$ws.Cells.Item(71, 14).Value = -16233.27285

# --- CUL!row 113 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 36296.066
$ws.Cells.Item(113, 10).Value = 847.0769
$ws.Cells.Item(113, 12).Value = 2541.2307
$ws.Cells.Item(113, 14).Value = -6881.2307

# --- CUL!row 131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 27874.934
$ws.Cells.Item(131, 10).Value = 28822.346
$ws.Cells.Item(131, 12).Value = 86467.038
$ws.Cells.Item(131, 14).Value = -96547.038

# --- CUL!row 135 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 560.9231
$ws.Cells.Item(135, 9).Value = 510.45
$ws.Cells.Item(135, 10).Value = 729.1667
$ws.Cells.Item(135, 11).Value = 4594.05
$ws.Cells.Item(135, 12).Value = 6562.5003
$ws.Cells.Item(135, 13).Value = -2059.05
$ws.Cells.Item(135, 14).Value = -11632.5003

# --- GSM!row 64 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).ClearContents()

# --- GSM!row 67 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).ClearContents()

# --- GSM!row 97 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 913.2222
$ws.Cells.Item(97, 9).Value = 940.1539
$ws.Cells.Item(97, 10).Value = 843.2
$ws.Cells.Item(97, 11).Value = 940.1539
$ws.Cells.Item(97, 12).Value = 843.2
$ws.Cells.Item(97, 13).Value = -444.1539
$ws.Cells.Item(97, 14).Value = -1835.2

# --- GSM!row 113 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1605.5454
$ws.Cells.Item(113, 9).Value = 1215.2
$ws.Cells.Item(113, 11).Value = 1215.2
$ws.Cells.Item(113, 13).Value = 954.8

# --- GSM!row 122 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1192.75
$ws.Cells.Item(122, 9).Value = 1136.1666
$ws.Cells.Item(122, 10).Value = 1277.625
$ws.Cells.Item(122, 11).Value = 3408.4998
$ws.Cells.Item(122, 12).Value = 3832.875
$ws.Cells.Item(122, 13).Value = -958.4998000000001
$ws.Cells.Item(122, 14).Value = -8732.875

# --- GSM!row 132 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1329130
$ws.Cells.Item(132, 9).Value = 1925636.8
$ws.Cells.Item(132, 10).Value = 3559.4443
$ws.Cells.Item(132, 11).Value = 5776910.4
$ws.Cells.Item(132, 12).Value = 10678.3329
$ws.Cells.Item(132, 13).Value = -5774380.4
$ws.Cells.Item(132, 14).Value = -15738.3329

# --- LTW!row 61 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 2621.5
$ws.Cells.Item(61, 10).Value = 2960
$ws.Cells.Item(61, 12).Value = 2960
$ws.Cells.Item(61, 14).Value = -3364

# --- LTW!row 100 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 1400
$ws.Cells.Item(100, 9).Value = 1400
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 1400
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).Value = -859
$ws.Cells.Item(100, 14).ClearContents()

# --- LTW!row 113 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 2621.5
$ws.Cells.Item(113, 10).Value = 2960
$ws.Cells.Item(113, 12).Value = 2960
$ws.Cells.Item(113, 14).Value = -7300

# --- LTW!row 132 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2288.4565
$ws.Cells.Item(132, 9).Value = 1625.7273
$ws.Cells.Item(132, 10).Value = 2895.9583
$ws.Cells.Item(132, 11).Value = 4877.1819
$ws.Cells.Item(132, 12).Value = 8687.874899999999
$ws.Cells.Item(132, 13).Value = -2347.1819
$ws.Cells.Item(132, 14).Value = -13747.8749

# --- WVR!row 96 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 14299.833
$ws.Cells.Item(96, 9).Value = 1000
$ws.Cells.Item(96, 11).Value = 1000
$ws.Cells.Item(96, 13).Value = 373

# --- WVR!row 113 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(113, 8).Value = 433.7619
$ws.Cells.Item(113, 9).Value = 297.22223
$ws.Cells.Item(113, 11).Value = 891.66669
$ws.Cells.Item(113, 13).Value = 1278.33331

# --- WVR!row 122 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 31156.186
$ws.Cells.Item(122, 9).Value = 41365.35
$ws.Cells.Item(122, 11).Value = 124096.05
$ws.Cells.Item(122, 13).Value = -121646.05

# --- WVR!row 126 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 2202.6316
$ws.Cells.Item(126, 9).Value = 2029.6666
$ws.Cells.Item(126, 10).Value = 2499.1428
$ws.Cells.Item(126, 11).Value = 6088.9998
$ws.Cells.Item(126, 12).Value = 7497.428400000001
$ws.Cells.Item(126, 13).Value = -3618.9998

Write-Host "Applied 252 value updates and 5 cell clears."
